$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Step 1 --------------------------------------------------------------
# Split the run "... hoje está ser" into several runs, wrapping "está" with
# gramStart/gramEnd proofErr markers, as a (fake) grammar-check pass would.
$rng = $d.Content
$found = $rng.Find.Execute(
    " foi porque ainda me doem memorais de tão aborrecido que hoje está ser",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $xmlSplit = '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t xml:space="preserve"> foi porque ainda me doem memorais de tão aborrecido que hoje </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>está</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> ser</w:t></w:r>' +
        '</w:p>'
    $rng.InsertXML($xmlSplit)
}

# --- Step 2 --------------------------------------------------------------
# The existing "_GoBack" bookmark sits right after "... ser" (before the
# trailing ", " run). It needs to move to the very end of the document
# (after the new last paragraph), so drop it here and re-add it later.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 3 --------------------------------------------------------------
# Append five new paragraphs after the first paragraph: two "spell-checked"
# all-caps/ lower Z runs, two empty paragraphs, and a final lowercase run.
$endRng = $d.Content
$endRng.Collapse(0)

$xmlNewParas =
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Zzzzzzzzzzzzzzzzzzzzzz</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Zzzzzzzzzzzzzz</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '"/>' +
    '<w:p xmlns:w="' + $wNs + '"/>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t>zzzzzzzzzzzzzzzzzzzz</w:t></w:r>' +
    '</w:p>'

$endRng.InsertXML($xmlNewParas)

# --- Step 4 --------------------------------------------------------------
# Re-create the "_GoBack" bookmark at the very end of the document (end of
# the last new paragraph, right after "zzzzzzzzzzzzzzzzzzzz").
#
# Adding a bookmark collapsed exactly at the absolute end of the document's
# content confuses this host's Bookmarks.Add (it silently resets the
# bookmark to start at 0). Work around it by appending a throw-away
# placeholder character, anchoring the bookmark just before it, and then
# deleting the placeholder again - the bookmark correctly collapses back
# to the (now real) end of the document.
$tailRng = $d.Content
$tailRng.Collapse(0)
$tailRng.InsertAfter("X")

$anchorStart = $d.Content.End - 1
$lastRng = $d.Range($anchorStart, $anchorStart)
$d.Bookmarks.Add("_GoBack", $lastRng)

$placeholderRng = $d.Range($anchorStart, $anchorStart + 1)
$placeholderRng.Text = ""
